$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.834.72"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.561.45"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.40"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.480"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.60"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.782.86"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.561.81"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.514"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.838.46"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.32"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.24"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.47"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.60"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.05"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0466"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.381.60"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.526"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.53"
$ws.Range("E43").Value = "  +4.69%  "
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.19"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.58"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.696.89"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.55"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  +3.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0978"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("E51").Value = "  +0.71%  "
